$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.843.79'
$ws.Range('E2').Value = '  +2.74%  '
$ws.Range('D3').Value = '3.443.95'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''579.07'
$ws.Range('E5').Value = '  +4.08%  '
$ws.Range('D6').Value = '''186.50'
$ws.Range('E6').Value = '  +6.60%  '
$ws.Range('D7').Value = '''0.630'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '3.436.59'
$ws.Range('E8').Value = '  +2.35%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').Value = '''0.171'
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('D11').Value = '''0.644'
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('D12').Value = '''56.28'
$ws.Range('E12').Value = '  +4.34%  '
$ws.Range('D13').Value = '''0.0000276'
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('D14').Value = '''9.39'
$ws.Range('E14').Value = '  +2.56%  '
$ws.Range('D15').Value = '3.998.97'
$ws.Range('E15').Value = '  +2.70%  '
$ws.Range('D16').Value = '''18.67'
$ws.Range('E16').Value = '  +2.33%  '
$ws.Range('D17').Value = '3.460.66'
$ws.Range('E17').Value = '  +3.41%  '
$ws.Range('D18').Value = '66.935.11'
$ws.Range('E18').Value = '  +3.20%  '
$ws.Range('D19').Value = '''12.06'
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('E20').Value = '  -2.23%  '
$ws.Range('D21').Value = '''1.02'
$ws.Range('E21').Value = '  +2.75%  '
$ws.Range('D22').Value = '''484.82'
$ws.Range('E22').Value = '  +7.62%  '
$ws.Range('D23').Value = '''5.29'
$ws.Range('E23').Value = '  +7.68%  '
$ws.Range('D24').Value = '''16.87'
$ws.Range('E24').Value = '  +22.94%  '
$ws.Range('D25').Value = '''4.34'
$ws.Range('E25').Value = '  +6.42%  '
$ws.Range('D26').Value = '''89.35'
$ws.Range('E26').Value = '  +3.05%  '
$ws.Range('D27').Value = '''2.95'
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('D28').Value = '''10.93'
$ws.Range('E28').Value = '  +1.50%  '
$ws.Range('D29').Value = '''9.02'
$ws.Range('E29').Value = '  +3.96%  '
$ws.Range('D30').Value = '''31.23'
$ws.Range('E30').Value = '  +0.66%  '
$ws.Range('D31').Value = '''7.23'
$ws.Range('E31').Value = '  +9.81%  '
$ws.Range('D32').Value = '''596.77'
$ws.Range('E32').Value = '  +3.64%  '
$ws.Range('D33').Value = '''11.73'
$ws.Range('E33').Value = '  +2.40%  '
$ws.Range('D34').Value = '''63.40'
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('D35').Value = '''0.111'
$ws.Range('E35').Value = '  +3.52%  '
$ws.Range('D36').Value = '''0.150'
$ws.Range('E36').Value = '  +5.73%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').Value = '''36.69'
$ws.Range('E38').Value = '  +2.75%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').Value = '''0.386'
$ws.Range('E39').Value = '  +3.82%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '''3.53'
$ws.Range('E40').Value = '  -2.60%  '
$ws.Range('D41').Value = '3.252.84'
$ws.Range('E41').Value = '  +5.46%  '
$ws.Range('D42').Value = '0.0₃0749'
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('D43').Value = '''2.90'
$ws.Range('E43').Value = '  +4.51%  '
$ws.Range('D44').Value = '''0.0430'
$ws.Range('E44').Value = '  +2.97%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '''2.79'
$ws.Range('E45').Value = '  +22.68%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').Value = '''2.51'
$ws.Range('E46').Value = '  +2.41%  '
$ws.Range('D47').Value = '''3.25'
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('D49').Value = '''3.30'
$ws.Range('E49').Value = '  +13.94%  '
$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D50').Value = '''1.00'
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '''8.69'
$ws.Range('E51').Value = '  +4.89%  '
